# Added shipping, discount and due date
# - Adds a new "Shipping Charges" column (L) with per-line values.
# - Fixes a handful of mis-copied "Shipping Address" (column K) cells so
#   every line that belongs to the same invoice shows the same address.
# - Updates the sheet's active selection to match the edited area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column L: "Shipping Charges" -------------------------------------
$ws.Range("L1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (match K1's header style)
$excel.CutCopyMode = $false
$ws.Range("L1").Value2 = "Shipping Charges"

$shippingCharges = @{
    2  = 100
    3  = 50
    4  = 0
    5  = 90
    6  = 120
    7  = 100
    8  = 70
    9  = 70
    10 = 70
    11 = 70
    12 = 100
    13 = 100
    14 = 200
    15 = 100
    16 = 100
    17 = 100
    18 = 70
    19 = 250
    20 = 250
    21 = 0
    22 = 0
    23 = 10
    24 = 10
    25 = 75
    26 = 75
    27 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 100
    32 = 90
    33 = 90
    34 = 90
    35 = 250
    36 = 80
    37 = 65
    38 = 0
}

foreach ($row in $shippingCharges.Keys) {
    $ws.Range("L$row").Value2 = $shippingCharges[$row]
}

# --- Fix column K (Shipping Address) for a few line items ------------------
# These lines belong to the same invoice as the row referenced and should
# carry the same shipping address instead of the stray value they had.
$ws.Range("K26").Value2 = $ws.Range("K25").Value2
$ws.Range("K35").Value2 = $ws.Range("K33").Value2
$ws.Range("K36").Value2 = $ws.Range("K33").Value2
$ws.Range("K37").Value2 = $ws.Range("K33").Value2
$ws.Range("K38").Value2 = $ws.Range("K33").Value2

# --- Update the on-screen selection -----------------------------------------
$ws.Range("K40").Select()
